$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Range("A38").Value2 = " Straight Kot Pantolon  Koyu Mavi "
$ws.Range("B38").Value2 = "350 Tl"
$ws.Range("C38").Value2 = "Jeans"
$ws.Range("D38").Value2 = "STRAİGHTKOYU.jpg"
$ws.Range("E38").Value2 = "%100 pamuklu kumaşı sayesinde gün boyu konfor sunar ve cildin nefes almasını sağlar.Normal bel kesimi ile rahat hareket etmenize olanak tanırken, düz paça tasarımı modern bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F38").Value2 = "Var"
$ws.Range("E38").Font.Name = "Arial"
$ws.Range("E38").Font.Color = 6710886

# Row 39
$ws.Range("A39").Value2 = " Straight Kot Pantolon  Açık  Mavi "
$ws.Range("B39").Value2 = "350 Tl"
$ws.Range("C39").Value2 = "Jeans"
$ws.Range("D39").Value2 = "STRAİGHTAÇIK.jpg"
$ws.Range("E39").Value2 = "%100 pamuklu kumaşı sayesinde gün boyu konfor sunar ve cildin nefes almasını sağlar.Normal bel kesimi ile rahat hareket etmenize olanak tanırken, düz paça tasarımı modern bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F39").Value2 = "Var"
$ws.Range("E39").Font.Name = "Arial"
$ws.Range("E39").Font.Color = 6710886

# Row 40
$ws.Range("A40").Value2 = "Regular Fit Kot Pantolon Seapoint"
$ws.Range("B40").Value2 = "450 Tl"
$ws.Range("C40").Value2 = "Jeans"
$ws.Range("D40").Value2 = "SEAPOİNT.jpg"
$ws.Range("E40").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F40").Value2 = "Var"
$ws.Range("E40").Font.Name = "Arial"
$ws.Range("E40").Font.Color = 6710886

# Row 41
$ws.Range("A41").Value2 = "Regular Fit Kot Pantolon Ocean"
$ws.Range("B41").Value2 = "450 Tl"
$ws.Range("C41").Value2 = "Jeans"
$ws.Range("D41").Value2 = "OCEAN.jpg"
$ws.Range("E41").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F41").Value2 = "Var"

# Row 42
$ws.Range("A42").Value2 = "Regular Fit Kot Pantolon Ren"
$ws.Range("B42").Value2 = "450 Tl"
$ws.Range("C42").Value2 = "Jeans"
$ws.Range("D42").Value2 = "REN.jpg"
$ws.Range("E42").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F42").Value2 = "Var"

# Row 43
$ws.Range("A43").Value2 = "Regular Fit Kot Pantolon Heinkel"
$ws.Range("B43").Value2 = "450 Tl"
$ws.Range("C43").Value2 = "Jeans"
$ws.Range("D43").Value2 = "HEİNKEL.jpg"
$ws.Range("E43").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F43").Value2 = "Var"

# Row 44
$ws.Range("A44").Value2 = "Regular Fit Kot Pantolon Hein"
$ws.Range("B44").Value2 = "450 Tl"
$ws.Range("C44").Value2 = "Jeans"
$ws.Range("D44").Value2 = "HEİN.jpg"
$ws.Range("E44").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F44").Value2 = "Var"

# Row 45
$ws.Range("A45").Value2 = "Regular Fit Kot Pantolon Hawker"
$ws.Range("B45").Value2 = "450 Tl"
$ws.Range("C45").Value2 = "Jeans"
$ws.Range("D45").Value2 = "HAWKER.jpg"
$ws.Range("E45").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F45").Value2 = "Var"

# Row 46
$ws.Range("A46").Value2 = "Regular Fit Kot Pantolon Forius"
$ws.Range("B46").Value2 = "450 Tl"
$ws.Range("C46").Value2 = "Jeans"
$ws.Range("D46").Value2 = "FORİUS.jpg"
$ws.Range("E46").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F46").Value2 = "Var"

# Row 47
$ws.Range("A47").Value2 = "Regular Fit Kot Pantolon Douglas"
$ws.Range("B47").Value2 = "450 Tl"
$ws.Range("C47").Value2 = "Jeans"
$ws.Range("D47").Value2 = "DOUGLAS.jpg"
$ws.Range("E47").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F47").Value2 = "Var"

# Row 48
$ws.Range("A48").Value2 = "Regular Fit Kot Pantolon Angry"
$ws.Range("B48").Value2 = "450 Tl"
$ws.Range("C48").Value2 = "Jeans"
$ws.Range("D48").Value2 = "ANGRY.jpg"
$ws.Range("E48").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."
$ws.Range("F48").Value2 = "Var"

# Row 49
$ws.Range("A49").Value2 = "Regular Fit Kot Kanvas Pantolon Taş"
$ws.Range("B49").Value2 = "450 Tl"
$ws.Range("C49").Value2 = "Jeans"
$ws.Range("D49").Value2 = "KANVATAŞ.jpg"
$ws.Range("E49").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.5 cepli tasarımı ve Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.29-30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

# Row 50
$ws.Range("A50").Value2 = "Regular Fit Kot Kanvas Pantolon Kahverengi"
$ws.Range("B50").Value2 = "450 Tl"
$ws.Range("C50").Value2 = "Jeans"
$ws.Range("D50").Value2 = "KANVASKAHVE.jpg"
$ws.Range("E50").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.5 cepli tasarımı ve Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.29-30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

# Row 51
$ws.Range("A51").Value2 = "Regular Fit Kot Kanvas Pantolon Gri"
$ws.Range("B51").Value2 = "450 Tl"
$ws.Range("C51").Value2 = "Jeans"
$ws.Range("D51").Value2 = "KANVASGRİ.jpg"
$ws.Range("E51").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.5 cepli tasarımı ve Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.29-30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

# Row 52
$ws.Range("A52").Value2 = "Regular Fit Kot Kanvas Pantolon Bej"
$ws.Range("B52").Value2 = "450 Tl"
$ws.Range("C52").Value2 = "Jeans"
$ws.Range("D52").Value2 = "KANVASBEJ.jpg"
$ws.Range("E52").Value2 = "%98 pamuk ve %2 spandex karışımı materyali sayesinde konforlu bir deneyim sunar.5 cepli tasarımı ve Regular fit kesimi ile vücut hatlarınıza uyum sağlayarak şık bir görünüm kazandırır.29-30-31-32-33-34-36-38-40 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

# Remove now-unused trailing rows (53-54 dropped out of range)
$ws.Range("A53:F54").ClearContents()

# View state: selection + top-left scroll position
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A52").Select()
